# Natmi following Dr Hou advice
# Update LR-pair communication scores for Tgm2-Itgb3: switch from single-sample
# values to 3-sample-averaged values, and add the missing "sCs" sending-cluster
# rows so the sheet covers the full 4x4 (ECs/FAPs/M2/sCs) sending/target grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 16,20
$data[0,0] = 'ECs'
$data[0,1] = 'Tgm2'
$data[0,2] = 'Itgb3'
$data[0,3] = 'ECs'
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 79.61246300000001
$data[0,7] = 238.837389
$data[0,8] = 0.6728436998494041
$data[0,9] = 0.6728436998494042
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 2.481489333333333
$data[0,13] = 7.444467999999999
$data[0,14] = 0.2345069082418988
$data[0,15] = 0.2345069082418987
$data[0,16] = 197.5574777348946
$data[0,17] = 1778.017299614052
$data[0,18] = 0.1577864957817239
$data[0,19] = 0.1577864957817239

$data[1,0] = 'ECs'
$data[1,1] = 'Tgm2'
$data[1,2] = 'Itgb3'
$data[1,3] = 'FAPs'
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 79.61246300000001
$data[1,7] = 238.837389
$data[1,8] = 0.6728436998494041
$data[1,9] = 0.6728436998494042
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 7.245227
$data[1,13] = 21.735681
$data[1,14] = 0.6846919551326144
$data[1,15] = 0.6846919551326142
$data[1,16] = 576.810366464101
$data[1,17] = 5191.293298176909
$data[1,18] = 0.4606906683485504
$data[1,19] = 0.4606906683485504

$data[2,0] = 'ECs'
$data[2,1] = 'Tgm2'
$data[2,2] = 'Itgb3'
$data[2,3] = 'M2'
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 79.61246300000001
$data[2,7] = 238.837389
$data[2,8] = 0.6728436998494041
$data[2,9] = 0.6728436998494042
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.2001876666666667
$data[2,13] = 0.600563
$data[2,14] = 0.01891823194544989
$data[2,15] = 0.01891823194544989
$data[2,16] = 15.93743320555633
$data[2,17] = 143.436898850007
$data[2,18] = 0.01272901317678569
$data[2,19] = 0.01272901317678569

$data[3,0] = 'ECs'
$data[3,1] = 'Tgm2'
$data[3,2] = 'Itgb3'
$data[3,3] = 'sCs'
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 79.61246300000001
$data[3,7] = 238.837389
$data[3,8] = 0.6728436998494041
$data[3,9] = 0.6728436998494042
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.6548283333333333
$data[3,13] = 1.964485
$data[3,14] = 0.06188290468003712
$data[3,15] = 0.06188290468003711
$data[3,16] = 52.13249645885167
$data[3,17] = 469.192468129665
$data[3,18] = 0.04163752254234418
$data[3,19] = 0.04163752254234417

$data[4,0] = 'FAPs'
$data[4,1] = 'Tgm2'
$data[4,2] = 'Itgb3'
$data[4,3] = 'ECs'
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 3.815058666666667
$data[4,7] = 11.445176
$data[4,8] = 0.03224291890608301
$data[4,9] = 0.03224291890608302
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 2.481489333333333
$data[4,13] = 7.444467999999999
$data[4,14] = 0.2345069082418988
$data[4,15] = 0.2345069082418987
$data[4,16] = 9.46702738737422
$data[4,17] = 85.20324648636799
$data[4,18] = 0.007561187225359791
$data[4,19] = 0.007561187225359791

$data[5,0] = 'FAPs'
$data[5,1] = 'Tgm2'
$data[5,2] = 'Itgb3'
$data[5,3] = 'FAPs'
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 3.815058666666667
$data[5,7] = 11.445176
$data[5,8] = 0.03224291890608301
$data[5,9] = 0.03224291890608302
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 7.245227
$data[5,13] = 21.735681
$data[5,14] = 0.6846919551326144
$data[5,15] = 0.6846919551326142
$data[5,16] = 27.64096605831734
$data[5,17] = 248.768694524856
$data[5,18] = 0.02207646718498831
$data[5,19] = 0.02207646718498831

$data[6,0] = 'FAPs'
$data[6,1] = 'Tgm2'
$data[6,2] = 'Itgb3'
$data[6,3] = 'M2'
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 3.815058666666667
$data[6,7] = 11.445176
$data[6,8] = 0.03224291890608301
$data[6,9] = 0.03224291890608302
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.2001876666666667
$data[6,13] = 0.600563
$data[6,14] = 0.01891823194544989
$data[6,15] = 0.01891823194544989
$data[6,16] = 0.7637276926764445
$data[6,17] = 6.873549234087999
$data[6,18] = 0.0006099790184636099
$data[6,19] = 0.0006099790184636099

$data[7,0] = 'FAPs'
$data[7,1] = 'Tgm2'
$data[7,2] = 'Itgb3'
$data[7,3] = 'sCs'
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 3.815058666666667
$data[7,7] = 11.445176
$data[7,8] = 0.03224291890608301
$data[7,9] = 0.03224291890608302
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.6548283333333333
$data[7,13] = 1.964485
$data[7,14] = 0.06188290468003712
$data[7,15] = 0.06188290468003711
$data[7,16] = 2.498208508262223
$data[7,17] = 22.48387657436
$data[7,18] = 0.001995285477271302
$data[7,19] = 0.001995285477271302

$data[8,0] = 'M2'
$data[8,1] = 'Tgm2'
$data[8,2] = 'Itgb3'
$data[8,3] = 'ECs'
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 30.45313966666667
$data[8,7] = 91.359419
$data[8,8] = 0.2573743154429307
$data[8,9] = 0.2573743154429307
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 2.481489333333333
$data[8,13] = 7.444467999999999
$data[8,14] = 0.2345069082418988
$data[8,15] = 0.2345069082418987
$data[8,16] = 75.56914124934355
$data[8,17] = 680.122271244092
$data[8,18] = 0.06035605497539685
$data[8,19] = 0.06035605497539685

$data[9,0] = 'M2'
$data[9,1] = 'Tgm2'
$data[9,2] = 'Itgb3'
$data[9,3] = 'FAPs'
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 30.45313966666667
$data[9,7] = 91.359419
$data[9,8] = 0.2573743154429307
$data[9,9] = 0.2573743154429307
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 7.245227
$data[9,13] = 21.735681
$data[9,14] = 0.6846919551326144
$data[9,15] = 0.6846919551326142
$data[9,16] = 220.6399097477043
$data[9,17] = 1985.759187729339
$data[9,18] = 0.1762221232415384
$data[9,19] = 0.1762221232415384

$data[10,0] = 'M2'
$data[10,1] = 'Tgm2'
$data[10,2] = 'Itgb3'
$data[10,3] = 'M2'
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 30.45313966666667
$data[10,7] = 91.359419
$data[10,8] = 0.2573743154429307
$data[10,9] = 0.2573743154429307
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.2001876666666667
$data[10,13] = 0.600563
$data[10,14] = 0.01891823194544989
$data[10,15] = 0.01891823194544989
$data[10,16] = 6.096342972544111
$data[10,17] = 54.86708675289699
$data[10,18] = 0.004869066996350748
$data[10,19] = 0.004869066996350748

$data[11,0] = 'M2'
$data[11,1] = 'Tgm2'
$data[11,2] = 'Itgb3'
$data[11,3] = 'sCs'
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 30.45313966666667
$data[11,7] = 91.359419
$data[11,8] = 0.2573743154429307
$data[11,9] = 0.2573743154429307
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.6548283333333333
$data[11,13] = 1.964485
$data[11,14] = 0.06188290468003712
$data[11,15] = 0.06188290468003711
$data[11,16] = 19.94157869269056
$data[11,17] = 179.474208234215
$data[11,18] = 0.01592707022964468
$data[11,19] = 0.01592707022964468

$data[12,0] = 'sCs'
$data[12,1] = 'Tgm2'
$data[12,2] = 'Itgb3'
$data[12,3] = 'ECs'
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 4.441711333333334
$data[12,7] = 13.325134
$data[12,8] = 0.03753906580158222
$data[12,9] = 0.03753906580158223
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 2.481489333333333
$data[12,13] = 7.444467999999999
$data[12,14] = 0.2345069082418988
$data[12,15] = 0.2345069082418987
$data[12,16] = 11.02205929541244
$data[12,17] = 99.19853365871199
$data[12,18] = 0.008803170259418241
$data[12,19] = 0.008803170259418241

$data[13,0] = 'sCs'
$data[13,1] = 'Tgm2'
$data[13,2] = 'Itgb3'
$data[13,3] = 'FAPs'
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 4.441711333333334
$data[13,7] = 13.325134
$data[13,8] = 0.03753906580158222
$data[13,9] = 0.03753906580158223
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 7.245227
$data[13,13] = 21.735681
$data[13,14] = 0.6846919551326144
$data[13,15] = 0.6846919551326142
$data[13,16] = 32.18120687847267
$data[13,17] = 289.630861906254
$data[13,18] = 0.02570269635753719
$data[13,19] = 0.02570269635753719

$data[14,0] = 'sCs'
$data[14,1] = 'Tgm2'
$data[14,2] = 'Itgb3'
$data[14,3] = 'M2'
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 4.441711333333334
$data[14,7] = 13.325134
$data[14,8] = 0.03753906580158222
$data[14,9] = 0.03753906580158223
$data[14,10] = 2
$data[14,11] = 0.6666666666666666
$data[14,12] = 0.2001876666666667
$data[14,13] = 0.600563
$data[14,14] = 0.01891823194544989
$data[14,15] = 0.01891823194544989
$data[14,16] = 0.889175827826889
$data[14,17] = 8.002582450442
$data[14,18] = 0.0007101727538498384
$data[14,19] = 0.0007101727538498384

$data[15,0] = 'sCs'
$data[15,1] = 'Tgm2'
$data[15,2] = 'Itgb3'
$data[15,3] = 'sCs'
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 4.441711333333334
$data[15,7] = 13.325134
$data[15,8] = 0.03753906580158222
$data[15,9] = 0.03753906580158223
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 0.6548283333333333
$data[15,13] = 1.964485
$data[15,14] = 0.06188290468003712
$data[15,15] = 0.06188290468003711
$data[15,16] = 2.908558429554445
$data[15,17] = 26.17702586599
$data[15,18] = 0.002323026430776954
$data[15,19] = 0.002323026430776954

$ws.Range("A2:T17").Value = $data

